$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 20120551
$ws.Range("B2").Value = "Trần VĨnh Phúc"

$ws.Range("A3").Value = 20120500
$ws.Range("B3").Value = "Võ Đức Huy"

$ws.Range("C3").Select()
